{"js": "// Apply the set of text replacements described by the commit diff.\n// Each entry is a unique, exact \"before\" string (found exactly once in the\n// original document: the date line plus the 25 division-fact table cells)\n// mapped to its \"after\" replacement. We search for every \"before\" string\n// first (so all lookups are resolved against the *original* document text,\n// avoiding any accidental re-matching of text that a later replacement\n// might introduce), then perform all the replacements.\nconst pairs = [\n  [\"2025-03-25 Tuesday\", \"2025-03-26 Wednesday\"],\n  [\"15\u00f74=3, 3\", \"95\u00f73=31, 2\"],\n  [\"49\u00f76=8, 1\", \"73\u00f73=24, 1\"],\n  [\"17\u00f75=3, 2\", \"95\u00f74=23, 3\"],\n  [\"85\u00f77=12, 1\", \"67\u00f73=22, 1\"],\n  [\"92\u00f74=23, 0\", \"39\u00f78=4, 7\"],\n  [\"32\u00f79=3, 5\", \"29\u00f75=5, 4\"],\n  [\"98\u00f72=49, 0\", \"34\u00f74=8, 2\"],\n  [\"94\u00f72=47, 0\", \"81\u00f74=20, 1\"],\n  [\"87\u00f77=12, 3\", \"92\u00f74=23, 0\"],\n  [\"25\u00f77=3, 4\", \"49\u00f72=24, 1\"],\n  [\"92\u00f72=46, 0\", \"10\u00f74=2, 2\"],\n  [\"83\u00f78=10, 3\", \"73\u00f78=9, 1\"],\n  [\"63\u00f75=12, 3\", \"55\u00f75=11, 0\"],\n  [\"47\u00f72=23, 1\", \"34\u00f75=6, 4\"],\n  [\"98\u00f78=12, 2\", \"20\u00f78=2, 4\"],\n  [\"84\u00f72=42, 0\", \"96\u00f76=16, 0\"],\n  [\"77\u00f73=25, 2\", \"59\u00f75=11, 4\"],\n  [\"34\u00f78=4, 2\", \"92\u00f78=11, 4\"],\n  [\"91\u00f77=13, 0\", \"71\u00f78=8, 7\"],\n  [\"44\u00f79=4, 8\", \"29\u00f75=5, 4\"],\n  [\"53\u00f79=5, 8\", \"23\u00f77=3, 2\"],\n  [\"14\u00f74=3, 2\", \"49\u00f76=8, 1\"],\n  [\"22\u00f78=2, 6\", \"64\u00f76=10, 4\"],\n  [\"21\u00f72=10, 1\", \"82\u00f72=41, 0\"],\n  [\"55\u00f74=13, 3\", \"74\u00f78=9, 2\"],\n];\n\nconst body = context.document.body;\n\n// Kick off a search for each \"before\" string.\nconst searchResults = pairs.map(([before]) =>\n  body.search(before, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach(r => r.load(\"items\"));\nawait context.sync();\n\n// Replace each match with its corresponding \"after\" string.\nfor (let i = 0; i < pairs.length; i++) {\n  const [before, after] = pairs[i];\n  const items = searchResults[i].items;\n  if (items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${before}\"`);\n  }\n  items[0].insertText(after, \"Replace\");\n}\n\nawait context.sync();\n\n", "ps1": "# Apply the set of text replacements described by the commit diff.\n# Each 'before' string is unique and occurs exactly once in the original\n# document (the date line plus the 25 division-fact table cells), so a\n# plain Find/Replace (wdReplaceAll) against the whole document is safe.\n# The pairs are applied in the same top-to-bottom / left-to-right order\n# they appear in the document, which guarantees each search runs before\n# any later step could coincidentally write matching text into the body.\n\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"2025-03-25 Tuesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-03-26 Wednesday\", 2)\nif (-not $found) { throw \"Could not find text to replace: 2025-03-25 Tuesday\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"15\u00f74=3, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"95\u00f73=31, 2\", 2)\nif (-not $found) { throw \"Could not find text to replace: 15\u00f74=3, 3\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"49\u00f76=8, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"73\u00f73=24, 1\", 2)\nif (-not $found) { throw \"Could not find text to replace: 49\u00f76=8, 1\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"17\u00f75=3, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"95\u00f74=23, 3\", 2)\nif (-not $found) { throw \"Could not find text to replace: 17\u00f75=3, 2\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"85\u00f77=12, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"67\u00f73=22, 1\", 2)\nif (-not $found) { throw \"Could not find text to replace: 85\u00f77=12, 1\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"92\u00f74=23, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"39\u00f78=4, 7\", 2)\nif (-not $found) { throw \"Could not find text to replace: 92\u00f74=23, 0\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"32\u00f79=3, 5\", $false, $false, $false, $false, $false, $true, 1, $false, \"29\u00f75=5, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 32\u00f79=3, 5\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"98\u00f72=49, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"34\u00f74=8, 2\", 2)\nif (-not $found) { throw \"Could not find text to replace: 98\u00f72=49, 0\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"94\u00f72=47, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"81\u00f74=20, 1\", 2)\nif (-not $found) { throw \"Could not find text to replace: 94\u00f72=47, 0\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"87\u00f77=12, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"92\u00f74=23, 0\", 2)\nif (-not $found) { throw \"Could not find text to replace: 87\u00f77=12, 3\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"25\u00f77=3, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00f72=24, 1\", 2)\nif (-not $found) { throw \"Could not find text to replace: 25\u00f77=3, 4\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"92\u00f72=46, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"10\u00f74=2, 2\", 2)\nif (-not $found) { throw \"Could not find text to replace: 92\u00f72=46, 0\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"83\u00f78=10, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"73\u00f78=9, 1\", 2)\nif (-not $found) { throw \"Could not find text to replace: 83\u00f78=10, 3\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"63\u00f75=12, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"55\u00f75=11, 0\", 2)\nif (-not $found) { throw \"Could not find text to replace: 63\u00f75=12, 3\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"47\u00f72=23, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"34\u00f75=6, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 47\u00f72=23, 1\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"98\u00f78=12, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"20\u00f78=2, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 98\u00f78=12, 2\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"84\u00f72=42, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"96\u00f76=16, 0\", 2)\nif (-not $found) { throw \"Could not find text to replace: 84\u00f72=42, 0\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"77\u00f73=25, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"59\u00f75=11, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 77\u00f73=25, 2\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"34\u00f78=4, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"92\u00f78=11, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 34\u00f78=4, 2\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"91\u00f77=13, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"71\u00f78=8, 7\", 2)\nif (-not $found) { throw \"Could not find text to replace: 91\u00f77=13, 0\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"44\u00f79=4, 8\", $false, $false, $false, $false, $false, $true, 1, $false, \"29\u00f75=5, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 44\u00f79=4, 8\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"53\u00f79=5, 8\", $false, $false, $false, $false, $false, $true, 1, $false, \"23\u00f77=3, 2\", 2)\nif (-not $found) { throw \"Could not find text to replace: 53\u00f79=5, 8\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"14\u00f74=3, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00f76=8, 1\", 2)\nif (-not $found) { throw \"Could not find text to replace: 14\u00f74=3, 2\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"22\u00f78=2, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"64\u00f76=10, 4\", 2)\nif (-not $found) { throw \"Could not find text to replace: 22\u00f78=2, 6\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"21\u00f72=10, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"82\u00f72=41, 0\", 2)\nif (-not $found) { throw \"Could not find text to replace: 21\u00f72=10, 1\" }\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$found = $r.Find.Execute(\"55\u00f74=13, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"74\u00f78=9, 2\", 2)\nif (-not $found) { throw \"Could not find text to replace: 55\u00f74=13, 3\" }\n\n"}
